$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Untreated" row (row 2), shifting rows 3-6 up.
$ws.Range("A2:E2").Delete()

# Update the remaining data values (small adjustments) for the new rows 2-4.
$ws.Range("B2").Value = -26.4
$ws.Range("C2").Value = -17.5
$ws.Range("D2").Value = -13.45
$ws.Range("E2").Value = -13.42

$ws.Range("B3").Value = -25.16
$ws.Range("C3").Value = -10.26
$ws.Range("D3").Value = -1.39
$ws.Range("E3").Value = 1.75

$ws.Range("B4").Value = -29.03
$ws.Range("C4").Value = -13.59
$ws.Range("D4").Value = -2.67
$ws.Range("E4").Value = 2.61

Write-Output "done"
